# Apply the changes described in the commit to 13_LibFormula.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Library_Formula")

# --- Rename IND_x -> INDICATOR_x (column C, rows 5-24); row5 text is unchanged ---
# (write these first so the shared-string table keeps UNION_INDICATORS/INDICATOR_*
# ahead of the renamed UpdateAnalysisUnit_IT/updateOutputTable_IT entries, matching
# how the workbook was actually produced)
$ws.Cells.Item(5, 3).Value = "UNION_INDICATORS"

$indicatorRows = [ordered]@{
    6  = "INDICATOR_2"
    7  = "INDICATOR_6"
    8  = "INDICATOR_7"
    9  = "INDICATOR_10"
    10 = "INDICATOR_11"
    11 = "INDICATOR_12"
    12 = "INDICATOR_13"
    13 = "INDICATOR_15"
    14 = "INDICATOR_16"
    15 = "INDICATOR_17"
    16 = "INDICATOR_18"
    17 = "INDICATOR_22"
    18 = "INDICATOR_26"
    19 = "INDICATOR_27"
    20 = "INDICATOR_28"
    21 = "INDICATOR_29"
    22 = "INDICATOR_31"
    23 = "INDICATOR_34"
    24 = "INDICATOR_35"
}

foreach ($row in $indicatorRows.Keys) {
    $ws.Cells.Item($row, 3).Value = $indicatorRows[$row]
}

# C6's style changes to match the plain style used by rows 2-5 (cellXfs index 1)
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rename the CREATE/MODIFY action rows (RETAIL -> IT) ---
$ws.Cells.Item(2, 3).Value = "UpdateAnalysisUnit_IT"
$ws.Cells.Item(4, 3).Value = "updateOutputTable_IT"

# Row 2 gains a new value in column F ("String")
$ws.Cells.Item(2, 6).Value = "String"

# --- Update the selected/visible range on the sheet ---
$ws.Activate()
$ws.Range("E7").Select()
$excel.ActiveWindow.ScrollColumn = 2
